$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1190320826869504
$ws.Range("C2").Value = 10.34677158129881
$ws.Range("D2").Value = 22.3905356188092
$ws.Range("E2").Value = 91228006295.30009
$ws.Range("G2").Value = 91228006328.15643

# Row 3
$ws.Range("B3").Value = 0.003208871385164791
$ws.Range("C3").Value = 53694773.49795976
$ws.Range("D3").Value = 31965208863585940
$ws.Range("E3").Value = 91228006295.30009
$ws.Range("G3").Value = 31965300145287008
